# Atualização de bases das ligas, do dia: 08-04-2024 às 21:28
# Swap the home/away match rows that were re-ordered in the source feed.
# For each pair of rows below, the entire row's data (columns B..AC) is
# exchanged between the two rows, while column A (the positional "id"
# index) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(2, 3),
    @(4, 5),
    @(11, 13),
    @(14, 15),
    @(46, 47),
    @(80, 81)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rangeA = $ws.Range("B${r1}:AC${r1}")
    $rangeB = $ws.Range("B${r2}:AC${r2}")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}
